$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts "Data" and everything
# after it one column to the right) and add the new "Ano" header.
$ws.Columns("F:F").Insert()
$ws.Range("F1").Value = "Ano"

# Match the width of the neighbouring "Procedência" column (closest
# achievable value for this column's width).
$ws.Columns("F:F").ColumnWidth = 10.5

# Restore the active cell selection to L1 (the "Dimensão largura (cm)"
# header, which moved from K1 to L1 because of the inserted column).
[void]$ws.Range("L1").Select()
